$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2, 3) {
    $ws.Range("D$row").Value = 0.092
    $ws.Range("E$row").Value = 0.0837
    $ws.Range("F$row").Value = 0.0356
    $ws.Range("K$row").Value = 46.5
    $ws.Range("L$row").Value = 0.3496240601503759
    $ws.Range("M$row").Value = 29.8
    $ws.Range("N$row").Value = 0.04275466284074605
    $ws.Range("O$row").Value = 0.6408602150537634
    $ws.Range("P$row").Value = 29.8
    $ws.Range("Q$row").Value = 0.04275466284074605
    $ws.Range("R$row").Value = 0.6408602150537634
    $ws.Range("U$row").Value = 91.5
    $ws.Range("V$row").Value = 0.1312769010043042
    $ws.Range("W$row").Value = 0.1025132275132275
    $ws.Range("X$row").Value = 0.05161894062365116
    $ws.Range("Y$row").Value = 0.05089428688957635
    $ws.Range("Z$row").Value = 0.1218841642228739
    $ws.Range("AB$row").Value = 0.03093211379054739
    $ws.Range("AC$row").Value = -0.03093211379054739
    $ws.Range("AD$row").Value = 842.8
    $ws.Range("AF$row").Value = 842.8
    $ws.Range("AG$row").Value = 751.3
    $ws.Range("AH$row").Value = 0.5473438108845304
    $ws.Range("AI$row").Value = 0.6509616127288176
    $ws.Range("AJ$row").Value = 0.5187461161361596
    $ws.Range("AK$row").Value = 0.6244182180851064
}
